$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3; this shifts the former rows 3-8 down to rows 4-9,
# matching the row the new "Bochum vs Alemannia Aachen" fixture was added at.
$ws.Rows(3).Insert()

# Helper: write a literal string into a cell without letting the automatic
# type inference reinterpret values such as "2025-10-09" as a date serial number.
# Temporarily forcing a text number format makes Excel accept the value as text;
# the format is then reset back to General/Normal so no extra cell formatting
# is left behind.
function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

# Row 2: Friendly Matches | Winterthur vs FC Vaduz
$ws.Range("A2").Value = "Friendly Matches"
Set-TextValue "B2" "2025-10-09"
$ws.Range("C2").Value = "08:30:00"
$ws.Range("D2").Value = "Winterthur"
$ws.Range("E2").Value = "FC Vaduz"
$ws.Range("F2").Value = 1.04
$ws.Range("G2").Value = 1000
$ws.Range("H2").Value = 1.04
$ws.Range("I2").Value = 1000
$ws.Range("J2").Value = 1.03
$ws.Range("K2").Value = 1000
$ws.Range("L2").Value = 1.01
$ws.Range("M2").Value = 1.01
$ws.Range("N2").Value = 1.3
$ws.Range("O2").Value = 1.14
$ws.Range("P2").Value = 1.3
$ws.Range("Q2").Value = 1.2
$ws.Range("R2").Value = 1.23
$ws.Range("S2").Value = 1.5
$ws.Range("T2").Value = 1.03
$ws.Range("U2").Value = 1.03
$ws.Range("V2").Value = 1.01
$ws.Range("W2").Value = 1.01
$ws.Range("X2").Value = 1000
$ws.Range("Y2").Value = 1000
$ws.Range("Z2").Value = 1000
$ws.Range("AA2").Value = 1000
$ws.Range("AB2").Value = 1000
$ws.Range("AC2").Value = 1000
$ws.Range("AD2").Value = 1000
$ws.Range("AE2").Value = 1000
$ws.Range("AF2").Value = 1000
$ws.Range("AG2").Value = 1000
$ws.Range("AH2").Value = 1000
$ws.Range("AI2").Value = 1000
$ws.Range("AJ2").Value = 1000
$ws.Range("AK2").Value = 1000
$ws.Range("AL2").Value = 1000
$ws.Range("AM2").Value = 1000
$ws.Range("AN2").Value = 1000
$ws.Range("AO2").Value = 1000

# Row 3: Friendly Matches | Bochum vs Alemannia Aachen
$ws.Range("A3").Value = "Friendly Matches"
Set-TextValue "B3" "2025-10-09"
$ws.Range("C3").Value = "09:00:00"
$ws.Range("D3").Value = "Bochum"
$ws.Range("E3").Value = "Alemannia Aachen"
$ws.Range("F3").Value = 1.04
$ws.Range("G3").Value = 1000
$ws.Range("H3").Value = 1.04
$ws.Range("I3").Value = 1000
$ws.Range("J3").Value = 1.01
$ws.Range("K3").Value = 950
$ws.Range("L3").Value = 1.01
$ws.Range("M3").Value = 1.01
$ws.Range("N3").Value = 1.28
$ws.Range("O3").Value = 1.14
$ws.Range("P3").Value = 1.28
$ws.Range("Q3").Value = 1.2
$ws.Range("R3").Value = 1.18
$ws.Range("S3").Value = 1.5
$ws.Range("T3").Value = 1.01
$ws.Range("U3").Value = 1.01
$ws.Range("V3").Value = 1.01
$ws.Range("W3").Value = 1.01
$ws.Range("X3").Value = 1000
$ws.Range("Y3").Value = 1000
$ws.Range("Z3").Value = 1000
$ws.Range("AA3").Value = 1000
$ws.Range("AB3").Value = 1000
$ws.Range("AC3").Value = 1000
$ws.Range("AD3").Value = 1000
$ws.Range("AE3").Value = 1000
$ws.Range("AF3").Value = 1000
$ws.Range("AG3").Value = 1000
$ws.Range("AH3").Value = 1000
$ws.Range("AI3").Value = 1000
$ws.Range("AJ3").Value = 1000
$ws.Range("AK3").Value = 1000
$ws.Range("AL3").Value = 1000
$ws.Range("AM3").Value = 1000
$ws.Range("AN3").Value = 1000
$ws.Range("AO3").Value = 1000

# Row 4: Brazilian Serie B | Ferroviaria vs Chapecoense
$ws.Range("A4").Value = "Brazilian Serie B"
Set-TextValue "B4" "2025-10-09"
$ws.Range("C4").Value = "19:00:00"
$ws.Range("D4").Value = "Ferroviaria"
$ws.Range("E4").Value = "Chapecoense"
$ws.Range("F4").Value = 2.62
$ws.Range("G4").Value = 2.76
$ws.Range("H4").Value = 3
$ws.Range("I4").Value = 3.15
$ws.Range("J4").Value = 3.3
$ws.Range("K4").Value = 3.4
$ws.Range("L4").Value = 1.49
$ws.Range("M4").Value = 1.09
$ws.Range("N4").Value = 3
$ws.Range("O4").Value = 1.41
$ws.Range("P4").Value = 1.7
$ws.Range("Q4").Value = 2.34
$ws.Range("R4").Value = 1.25
$ws.Range("S4").Value = 4.3
$ws.Range("T4").Value = 1.91
$ws.Range("U4").Value = 1.96
$ws.Range("V4").Value = 1.46
$ws.Range("W4").Value = 1.57
$ws.Range("X4").Value = 11.5
$ws.Range("Y4").Value = 10.5
$ws.Range("Z4").Value = 20
$ws.Range("AA4").Value = 55
$ws.Range("AB4").Value = 10.5
$ws.Range("AC4").Value = 7.8
$ws.Range("AD4").Value = 13.5
$ws.Range("AE4").Value = 40
$ws.Range("AF4").Value = 18
$ws.Range("AG4").Value = 13
$ws.Range("AH4").Value = 21
$ws.Range("AI4").Value = 55
$ws.Range("AJ4").Value = 44
$ws.Range("AK4").Value = 36
$ws.Range("AL4").Value = 55
$ws.Range("AM4").Value = 150
$ws.Range("AN4").Value = 1000
$ws.Range("AO4").Value = 1000

# Row 5: Brazilian Serie B | Coritiba vs Atletico GO
$ws.Range("A5").Value = "Brazilian Serie B"
Set-TextValue "B5" "2025-10-09"
$ws.Range("C5").Value = "19:30:00"
$ws.Range("D5").Value = "Coritiba"
$ws.Range("E5").Value = "Atletico GO"
$ws.Range("F5").Value = 1.87
$ws.Range("G5").Value = 1.97
$ws.Range("H5").Value = 5.3
$ws.Range("I5").Value = 6
$ws.Range("J5").Value = 3.3
$ws.Range("K5").Value = 3.45
$ws.Range("L5").Value = 1.58
$ws.Range("M5").Value = 1.12
$ws.Range("N5").Value = 2.54
$ws.Range("O5").Value = 1.56
$ws.Range("P5").Value = 1.51
$ws.Range("Q5").Value = 2.66
$ws.Range("R5").Value = 1.17
$ws.Range("S5").Value = 5.7
$ws.Range("T5").Value = 2.28
$ws.Range("U5").Value = 1.64
$ws.Range("V5").Value = 1.2
$ws.Range("W5").Value = 2.04
$ws.Range("X5").Value = 9.800000000000001
$ws.Range("Y5").Value = 980
$ws.Range("Z5").Value = 42
$ws.Range("AA5").Value = 200
$ws.Range("AB5").Value = 6.4
$ws.Range("AC5").Value = 8
$ws.Range("AD5").Value = 980
$ws.Range("AE5").Value = 130
$ws.Range("AF5").Value = 970
$ws.Range("AG5").Value = 980
$ws.Range("AH5").Value = 980
$ws.Range("AI5").Value = 140
$ws.Range("AJ5").Value = 24
$ws.Range("AK5").Value = 34
$ws.Range("AL5").Value = 70
$ws.Range("AM5").Value = 260
$ws.Range("AN5").Value = 26
$ws.Range("AO5").Value = 240

# Row 6: Chilean Primera Division | Univ Catolica (Chile) vs Nublense
$ws.Range("A6").Value = "Chilean Primera Division"
Set-TextValue "B6" "2025-10-09"
$ws.Range("C6").Value = "20:30:00"
$ws.Range("D6").Value = "Univ Catolica (Chile)"
$ws.Range("E6").Value = "Nublense"
$ws.Range("F6").Value = 1.87
$ws.Range("G6").Value = 1.94
$ws.Range("H6").Value = 4.8
$ws.Range("I6").Value = 5.4
$ws.Range("J6").Value = 3.5
$ws.Range("K6").Value = 3.8
$ws.Range("L6").Value = 1.42
$ws.Range("M6").Value = 1.09
$ws.Range("N6").Value = 3.1
$ws.Range("O6").Value = 1.41
$ws.Range("P6").Value = 1.73
$ws.Range("Q6").Value = 2.2
$ws.Range("R6").Value = 1.27
$ws.Range("S6").Value = 4.1
$ws.Range("T6").Value = 2.02
$ws.Range("U6").Value = 1.84
$ws.Range("V6").Value = 1.23
$ws.Range("W6").Value = 2.06
$ws.Range("X6").Value = 980
$ws.Range("Y6").Value = 980
$ws.Range("Z6").Value = 980
$ws.Range("AA6").Value = 130
$ws.Range("AB6").Value = 980
$ws.Range("AC6").Value = 980
$ws.Range("AD6").Value = 980
$ws.Range("AE6").Value = 80
$ws.Range("AF6").Value = 980
$ws.Range("AG6").Value = 980
$ws.Range("AH6").Value = 980
$ws.Range("AI6").Value = 100
$ws.Range("AJ6").Value = 980
$ws.Range("AK6").Value = 980
$ws.Range("AL6").Value = 980
$ws.Range("AM6").Value = 180
$ws.Range("AN6").Value = 980
$ws.Range("AO6").Value = 120

# Row 7: FIFA World Cup Qualifiers - Americas | Nicaragua vs Haiti
$ws.Range("A7").Value = "FIFA World Cup Qualifiers - Americas"
Set-TextValue "B7" "2025-10-09"
$ws.Range("C7").Value = "21:00:00"
$ws.Range("D7").Value = "Nicaragua"
$ws.Range("E7").Value = "Haiti"
$ws.Range("F7").Value = 3.15
$ws.Range("G7").Value = 3.5
$ws.Range("H7").Value = 2.36
$ws.Range("I7").Value = 2.48
$ws.Range("J7").Value = 3.45
$ws.Range("K7").Value = 3.85
$ws.Range("L7").Value = 1.43
$ws.Range("M7").Value = 1.07
$ws.Range("N7").Value = 3.55
$ws.Range("O7").Value = 1.33
$ws.Range("P7").Value = 1.88
$ws.Range("Q7").Value = 2.08
$ws.Range("R7").Value = 1.32
$ws.Range("S7").Value = 3.85
$ws.Range("T7").Value = 1.76
$ws.Range("U7").Value = 2.06
$ws.Range("V7").Value = 1.67
$ws.Range("W7").Value = 1.4
$ws.Range("X7").Value = 14.5
$ws.Range("Y7").Value = 11
$ws.Range("Z7").Value = 16
$ws.Range("AA7").Value = 34
$ws.Range("AB7").Value = 13
$ws.Range("AC7").Value = 8.6
$ws.Range("AD7").Value = 12
$ws.Range("AE7").Value = 27
$ws.Range("AF7").Value = 23
$ws.Range("AG7").Value = 15
$ws.Range("AH7").Value = 19.5
$ws.Range("AI7").Value = 980
$ws.Range("AJ7").Value = 980
$ws.Range("AK7").Value = 40
$ws.Range("AL7").Value = 55
$ws.Range("AM7").Value = 130
$ws.Range("AN7").Value = 38
$ws.Range("AO7").Value = 22

# Row 8: Brazilian Serie B | Remo vs Athletico-PR
$ws.Range("A8").Value = "Brazilian Serie B"
Set-TextValue "B8" "2025-10-09"
$ws.Range("C8").Value = "21:35:00"
$ws.Range("D8").Value = "Remo"
$ws.Range("E8").Value = "Athletico-PR"
$ws.Range("F8").Value = 3.45
$ws.Range("G8").Value = 3.75
$ws.Range("H8").Value = 2.32
$ws.Range("I8").Value = 2.48
$ws.Range("J8").Value = 3.1
$ws.Range("K8").Value = 3.4
$ws.Range("L8").Value = 1.56
$ws.Range("M8").Value = 1.11
$ws.Range("N8").Value = 2.78
$ws.Range("O8").Value = 1.47
$ws.Range("P8").Value = 1.6
$ws.Range("Q8").Value = 2.42
$ws.Range("R8").Value = 1.22
$ws.Range("S8").Value = 4.8
$ws.Range("T8").Value = 2.06
$ws.Range("U8").Value = 1.83
$ws.Range("V8").Value = 1.68
$ws.Range("W8").Value = 1.37
$ws.Range("X8").Value = 9.6
$ws.Range("Y8").Value = 8.4
$ws.Range("Z8").Value = 13.5
$ws.Range("AA8").Value = 34
$ws.Range("AB8").Value = 11
$ws.Range("AC8").Value = 7.4
$ws.Range("AD8").Value = 12
$ws.Range("AE8").Value = 32
$ws.Range("AF8").Value = 1000
$ws.Range("AG8").Value = 18
$ws.Range("AH8").Value = 1000
$ws.Range("AI8").Value = 60
$ws.Range("AJ8").Value = 95
$ws.Range("AK8").Value = 1000
$ws.Range("AL8").Value = 1000
$ws.Range("AM8").Value = 210
$ws.Range("AN8").Value = 1000
$ws.Range("AO8").Value = 1000

# Row 9: FIFA World Cup Qualifiers - Americas | Honduras vs Costa Rica
$ws.Range("A9").Value = "FIFA World Cup Qualifiers - Americas"
Set-TextValue "B9" "2025-10-09"
$ws.Range("C9").Value = "23:00:00"
$ws.Range("D9").Value = "Honduras"
$ws.Range("E9").Value = "Costa Rica"
$ws.Range("F9").Value = 2.8
$ws.Range("G9").Value = 980
$ws.Range("H9").Value = 1.04
$ws.Range("I9").Value = 980
$ws.Range("J9").Value = 1.2
$ws.Range("K9").Value = 980
$ws.Range("L9").Value = 1.46
$ws.Range("M9").Value = 1.06
$ws.Range("N9").Value = 1.11
$ws.Range("O9").Value = 1.36
$ws.Range("P9").Value = 1.24
$ws.Range("Q9").Value = 1.36
$ws.Range("R9").Value = 1.18
$ws.Range("S9").Value = 1.05
$ws.Range("T9").Value = 1.03
$ws.Range("U9").Value = 1.03
$ws.Range("V9").Value = 1.01
$ws.Range("W9").Value = 1.38
$ws.Range("X9").Value = 1000
$ws.Range("Y9").Value = 1000
$ws.Range("Z9").Value = 1000
$ws.Range("AA9").Value = 1000
$ws.Range("AB9").Value = 1000
$ws.Range("AC9").Value = 1000
$ws.Range("AD9").Value = 1000
$ws.Range("AE9").Value = 1000
$ws.Range("AF9").Value = 1000
$ws.Range("AG9").Value = 1000
$ws.Range("AH9").Value = 1000
$ws.Range("AI9").Value = 1000
$ws.Range("AJ9").Value = 1000
$ws.Range("AK9").Value = 1000
$ws.Range("AL9").Value = 1000
$ws.Range("AM9").Value = 1000
$ws.Range("AN9").Value = 1000
$ws.Range("AO9").Value = 1000
